# "Generate Report for Handoff": a new handoff/handback cycle produced a
# fresh localization id + content hash, so every cell (and the matching
# hyperlink display text) that referenced the previous id/hash/timestamps
# must be refreshed to the new ones.

$oldId = "58a7cc9d-2090-4c07-bda3-cbcd62a034ff"
$newId = "1e5cb89d-5659-48d6-9c61-93610b0ac39f"

$oldHash = "057508a69b3db97ff2e19d7ccbc6c6ed1635d349"
$newHash = "09aa81012602d30b190cf72cc7991add4a918cf5"

# Same external GitHub blob link is reused by every sheet's hyperlink (only
# the visible display text changes to match the new id).
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/21182ba393df0821caf60589e43c6d71edfb75ea/e2e/$oldId.md"

$wb = $excel.ActiveWorkbook

function Update-Hyperlink($ws, $cellRef, $displayText) {
    $rng = $ws.Range($cellRef)
    $rng.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($rng, $linkAddress, "", "", $displayText)
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-09-07 01:17:09"
Update-Hyperlink $wsOverview "B2" "e2e\$newId.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-07 01:16:58"
Update-Hyperlink $wsZh "A2" "$newId.md"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-07 01:17:09"
Update-Hyperlink $wsDe "A2" "$newId.md"
